# 996: Add T2A to xlsx files and extract process
#
# WMT_Extract's header row (row 1) is selected and copied, a new worksheet
# "T2A" is appended after the last sheet (GS), and the header row is pasted
# into it as the start of a fresh extract tab.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("WMT_Extract")
$src.Activate()
$headerRow = $src.Range("A1:XFD1")
$headerRow.Select()
$headerRow.Copy()

$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $wsLast)
$ws.Name = "T2A"

$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A1").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

$ws.Rows.Item(1).RowHeight = 17
$ws.Range("A1:XFD1").Select()
